$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.05619466666666667
$ws.Range("I2").Value = 0.04986276087265156
$ws.Range("J2").Value = 0.07297477932340853
$ws.Range("M2").Value = 0.8584576666666667
$ws.Range("N2").Value = 2.575373
$ws.Range("O2").Value = 0.02952026538348031
$ws.Range("P2").Value = 0.03028938521394646
$ws.Range("Q2").Value = 0.04824074242577778
$ws.Range("R2").Value = 0.434166681832
$ws.Range("S2").Value = 0.001471961933713692
$ws.Range("T2").Value = 0.002210361201829456

$ws.Range("G3").Value = 0.05619466666666667
$ws.Range("I3").Value = 0.04986276087265156
$ws.Range("J3").Value = 0.07297477932340853
$ws.Range("N3").Value = 75.717583
$ws.Range("O3").Value = 0.8679143348771993
$ws.Range("P3").Value = 0.8905269407406087
$ws.Range("S3").Value = 0.04327660493792822
$ws.Range("T3").Value = 0.06498600698209601

$ws.Range("G4").Value = 0.05619466666666667
$ws.Range("I4").Value = 0.04986276087265156
$ws.Range("J4").Value = 0.07297477932340853
$ws.Range("M4").Value = 0.3580240000000001
$ws.Range("N4").Value = 1.074072
$ws.Range("O4").Value = 0.0123115721415754
$ws.Range("P4").Value = 0.01263233735676886
$ws.Range("Q4").Value = 0.02011903933866667
$ws.Range("R4").Value = 0.181071354048
$ws.Range("S4").Value = 0.000613888977661773
$ws.Range("T4").Value = 0.0009218420309490576

$ws.Range("G5").Value = 0.05619466666666667
$ws.Range("I5").Value = 0.04986276087265156
$ws.Range("J5").Value = 0.07297477932340853
$ws.Range("M5").Value = 2.2152535
$ws.Range("N5").Value = 4.430507
$ws.Range("O5").Value = 0.07617716487477769
$ws.Range("P5").Value = 0.05210792115009603
$ws.Range("Q5").Value = 0.1244854320146667
$ws.Range("R5").Value = 0.7469125920880001
$ws.Range("S5").Value = 0.003798403756107592
$ws.Range("T5").Value = 0.00380256404692983

$ws.Range("G6").Value = 0.05619466666666667
$ws.Range("I6").Value = 0.04986276087265156
$ws.Range("J6").Value = 0.07297477932340853
$ws.Range("M6").Value = 0.4093533333333334
$ws.Range("N6").Value = 1.22806
$ws.Range("O6").Value = 0.01407666272296744
$ws.Range("P6").Value = 0.01444341553857988
$ws.Range("Q6").Value = 0.02300347411555556
$ws.Range("R6").Value = 0.20703126704
$ws.Range("S6").Value = 0.0007019012672402939
$ws.Range("T6").Value = 0.001054005061604157

$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.070792
$ws.Range("H7").Value = 2.141584
$ws.Range("I7").Value = 0.9501372391273485
$ws.Range("J7").Value = 0.9270252206765914
$ws.Range("M7").Value = 0.8584576666666667
$ws.Range("N7").Value = 2.575373
$ws.Range("O7").Value = 0.02952026538348031
$ws.Range("P7").Value = 0.03028938521394646
$ws.Range("Q7").Value = 0.9192296018053333
$ws.Range("R7").Value = 5.515377610832
$ws.Range("S7").Value = 0.02804830344976662
$ws.Range("T7").Value = 0.028079024012117

$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.070792
$ws.Range("H8").Value = 2.141584
$ws.Range("I8").Value = 0.9501372391273485
$ws.Range("J8").Value = 0.9270252206765914
$ws.Range("N8").Value = 75.717583
$ws.Range("O8").Value = 0.8679143348771993
$ws.Range("P8").Value = 0.8905269407406087
$ws.Range("Q8").Value = 27.02592737857866
$ws.Range("R8").Value = 162.155564271472
$ws.Range("S8").Value = 0.8246377299392711
$ws.Range("T8").Value = 0.8255409337585126

$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.070792
$ws.Range("H9").Value = 2.141584
$ws.Range("I9").Value = 0.9501372391273485
$ws.Range("J9").Value = 0.9270252206765914
$ws.Range("M9").Value = 0.3580240000000001
$ws.Range("N9").Value = 1.074072
$ws.Range("O9").Value = 0.0123115721415754
$ws.Range("P9").Value = 0.01263233735676886
$ws.Range("Q9").Value = 0.3833692350080001
$ws.Range("R9").Value = 2.300215410048
$ws.Range("S9").Value = 0.01169768316391363
$ws.Range("T9").Value = 0.01171049532581981

$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.070792
$ws.Range("H10").Value = 2.141584
$ws.Range("I10").Value = 0.9501372391273485
$ws.Range("J10").Value = 0.9270252206765914
$ws.Range("M10").Value = 2.2152535
$ws.Range("N10").Value = 4.430507
$ws.Range("O10").Value = 0.07617716487477769
$ws.Range("P10").Value = 0.05210792115009603
$ws.Range("Q10").Value = 2.372075725772
$ws.Range("R10").Value = 9.488302903088
$ws.Range("S10").Value = 0.0723787611186701
$ws.Range("T10").Value = 0.0483053571031662

$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.070792
$ws.Range("H11").Value = 2.141584
$ws.Range("I11").Value = 0.9501372391273485
$ws.Range("J11").Value = 0.9270252206765914
$ws.Range("M11").Value = 0.4093533333333334
$ws.Range("N11").Value = 1.22806
$ws.Range("O11").Value = 0.01407666272296744
$ws.Range("P11").Value = 0.01444341553857988
$ws.Range("Q11").Value = 0.4383322745066667
$ws.Range("R11").Value = 2.62999364704
$ws.Range("S11").Value = 0.01337476145572715
$ws.Range("T11").Value = 0.01338941047697573
